# Generate Report for Handoff
# Adds two new localization-pipeline entries (2344ce4e-... and 83f1369c-...)
# to the existing "Overview", "zh-cn" and "de-de" report sheets, inserted so
# that 2344ce4e-... lands in row 3 (pushing the former row-3 entry,
# 73f1d282-..., down to row 4) and 83f1369c-... is appended as row 5.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": columns A..G
#   A File Name, B Path And Name, C Extension, D Publish URL,
#   E zh-cn, F de-de, G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Shift the old row 3 (73f1d282) down to row 4, then write the new
# row 3 (2344ce4e) and new row 5 (83f1369c).
$wsOverview.Rows.Item(3).Copy()
$wsOverview.Rows.Item(4).PasteSpecial()

$wsOverview.Range("A3").Value2 = "2344ce4e-f818-412b-87b5-3219049dcd59.md"
$wsOverview.Range("B3").Value2 = "e2e\2344ce4e-f818-412b-87b5-3219049dcd59.md"
$wsOverview.Range("C3").Value2 = ".md"
$wsOverview.Range("D3").Value2 = ""
$wsOverview.Range("E3").Value2 = "Ready for handoff"
$wsOverview.Range("F3").Value2 = "Ready for handoff"
$wsOverview.Range("G3").Value2 = "2016-08-27 20:40:24"

$wsOverview.Range("A5").Value2 = "83f1369c-36ff-4152-b1b0-c9170391a0dc.md"
$wsOverview.Range("B5").Value2 = "e2e\83f1369c-36ff-4152-b1b0-c9170391a0dc.md"
$wsOverview.Range("C5").Value2 = ".md"
$wsOverview.Range("D5").Value2 = ""
$wsOverview.Range("E5").Value2 = "Ready for handoff"
$wsOverview.Range("F5").Value2 = "Ready for handoff"
$wsOverview.Range("G5").Value2 = "2016-08-27 20:40:24"
$wsOverview.Range("G5").NumberFormat = $wsOverview.Range("G4").NumberFormat

# Hyperlinks on column B (file name links) - rebuild to match the new
# row order: B2 e6f217b5 (unchanged), B3 2344ce4e, B4 73f1d282, B5 83f1369c
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/28a6af83fcb9da17550f5b8f644dc1a16401668e/e2e/e6f217b5-9e9e-4642-8dad-93684e3a2563.md", "", "", "e2e\e6f217b5-9e9e-4642-8dad-93684e3a2563.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2344ce4ef818412b87b53219049dcd59aaaaaaa/e2e/2344ce4e-f818-412b-87b5-3219049dcd59.md", "", "", "e2e\2344ce4e-f818-412b-87b5-3219049dcd59.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/70112b919a6e5dc5d2ab3920ed16a98ae3fc3b16/e2e/73f1d282-66db-476d-9a58-30ea6fde4634.md", "", "", "e2e\73f1d282-66db-476d-9a58-30ea6fde4634.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/83f1369c36ff4152b1b0c9170391a0dcaaaaaaa/e2e/83f1369c-36ff-4152-b1b0-c9170391a0dc.md", "", "", "e2e\83f1369c-36ff-4152-b1b0-c9170391a0dc.md")

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G5"))

# ---------------------------------------------------------------------
# Shared layout/content for the language detail sheets ("zh-cn", "de-de")
# Columns A..P:
#  A Source File Name, B File Extension, C Status, D Source Path,
#  E Priority, F Content Duplicate, G Latest Handoff File,
#  H Latest Handoff Datetime, I Latest Target File, J Latest Handback File,
#  K Latest Handback DateTime, L Reference Tokens, M To be localized,
#  N Dependency From, O Has metadata, P Error Detail
# ---------------------------------------------------------------------
function Set-LangSheetRow {
    param($ws, $row, $vals)
    $cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value2 = $vals[$i]
    }
}

# ---- zh-cn ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Rows.Item(3).Copy()
$wsZh.Rows.Item(4).PasteSpecial()

Set-LangSheetRow $wsZh 3 @(
    "2344ce4e-f818-412b-87b5-3219049dcd59.md", ".md", "Ready for handoff", "e2e", "ht", "False",
    "2344ce4e-f818-412b-87b5-3219049dcd59.af3057e2e383ada60d7b6006fd204039fb19be07.zh-cn.xlf",
    "2016-08-27 20:40:19", "", "", "0001-01-01 00:00:00", "", "True", "", "False", ""
)

Set-LangSheetRow $wsZh 5 @(
    "83f1369c-36ff-4152-b1b0-c9170391a0dc.md", ".md", "Ready for handoff", "e2e", "ht", "False",
    "83f1369c-36ff-4152-b1b0-c9170391a0dc.e41edec302dea881c6982bac0ffac5cc6a56d120.zh-cn.xlf",
    "2016-08-27 20:40:19", "", "", "0001-01-01 00:00:00", "", "True", "", "False", ""
)
$wsZh.Range("H5").NumberFormat = $wsZh.Range("H4").NumberFormat
$wsZh.Range("K5").NumberFormat = $wsZh.Range("K4").NumberFormat

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/28a6af83fcb9da17550f5b8f644dc1a16401668e/e2e/e6f217b5-9e9e-4642-8dad-93684e3a2563.md", "", "", "e6f217b5-9e9e-4642-8dad-93684e3a2563.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/94a687f56d44b6c705add0d877f97614be440927/e2e/e6f217b5-9e9e-4642-8dad-93684e3a2563.md", "", "", "e6f217b5-9e9e-4642-8dad-93684e3a2563.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2344ce4ef818412b87b53219049dcd59aaaaaaa/e2e/2344ce4e-f818-412b-87b5-3219049dcd59.md", "", "", "2344ce4e-f818-412b-87b5-3219049dcd59.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/70112b919a6e5dc5d2ab3920ed16a98ae3fc3b16/e2e/73f1d282-66db-476d-9a58-30ea6fde4634.md", "", "", "73f1d282-66db-476d-9a58-30ea6fde4634.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/83f1369c36ff4152b1b0c9170391a0dcaaaaaaa/e2e/83f1369c-36ff-4152-b1b0-c9170391a0dc.md", "", "", "83f1369c-36ff-4152-b1b0-c9170391a0dc.md")

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P5"))

# ---- de-de ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Rows.Item(3).Copy()
$wsDe.Rows.Item(4).PasteSpecial()

Set-LangSheetRow $wsDe 3 @(
    "2344ce4e-f818-412b-87b5-3219049dcd59.md", ".md", "Ready for handoff", "e2e", "ht", "False",
    "2344ce4e-f818-412b-87b5-3219049dcd59.af3057e2e383ada60d7b6006fd204039fb19be07.de-de.xlf",
    "2016-08-27 20:40:24", "", "", "0001-01-01 00:00:00", "", "True", "", "False", ""
)

Set-LangSheetRow $wsDe 5 @(
    "83f1369c-36ff-4152-b1b0-c9170391a0dc.md", ".md", "Ready for handoff", "e2e", "ht", "False",
    "83f1369c-36ff-4152-b1b0-c9170391a0dc.e41edec302dea881c6982bac0ffac5cc6a56d120.de-de.xlf",
    "2016-08-27 20:40:24", "", "", "0001-01-01 00:00:00", "", "True", "", "False", ""
)
$wsDe.Range("H5").NumberFormat = $wsDe.Range("H4").NumberFormat
$wsDe.Range("K5").NumberFormat = $wsDe.Range("K4").NumberFormat

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/28a6af83fcb9da17550f5b8f644dc1a16401668e/e2e/e6f217b5-9e9e-4642-8dad-93684e3a2563.md", "", "", "e6f217b5-9e9e-4642-8dad-93684e3a2563.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/cdd747fb76ca4ffa0bbf200c15aad201925f8231/e2e/e6f217b5-9e9e-4642-8dad-93684e3a2563.md", "", "", "e6f217b5-9e9e-4642-8dad-93684e3a2563.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2344ce4ef818412b87b53219049dcd59aaaaaaa/e2e/2344ce4e-f818-412b-87b5-3219049dcd59.md", "", "", "2344ce4e-f818-412b-87b5-3219049dcd59.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/70112b919a6e5dc5d2ab3920ed16a98ae3fc3b16/e2e/73f1d282-66db-476d-9a58-30ea6fde4634.md", "", "", "73f1d282-66db-476d-9a58-30ea6fde4634.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/83f1369c36ff4152b1b0c9170391a0dcaaaaaaa/e2e/83f1369c-36ff-4152-b1b0-c9170391a0dc.md", "", "", "83f1369c-36ff-4152-b1b0-c9170391a0dc.md")

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P5"))
